# New .ttl from Google sheet has been generated:
#  - refresh the "pav:lastUpdatedOn" timestamp (B16)
#  - append two new vocabulary term rows (21 and 22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "lastUpdatedOn" timestamp
$ws.Range("B16").Value = "2021-08-24T09:33+02:00"

# Row 21: beer-onto:beer_nutrition
$ws.Cells.Item(21, 1).Value = "beer-onto:beer_nutrition"
$ws.Cells.Item(21, 2).Value = "beer_nutrition"
$ws.Cells.Item(21, 4).Value = "Calories in a unit of a beer."

# Row 22: beer-onto:beer_bitterness
$ws.Cells.Item(22, 1).Value = "beer-onto:beer_bitterness"
$ws.Cells.Item(22, 2).Value = "beer_bitterness"
$ws.Cells.Item(22, 4).Value = "The scale of beer bitterness determining the taste of beer."
